# Generate Report for Handoff
# Replaces the old localization GUID/file references with the new ones,
# and refreshes the handoff/handback timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "c52c2fd3-8ba9-43c6-9ddf-68357de1540e"
$newGuid = "d1d7db9c-142b-4746-a084-33bc9f3a0d22"
$oldHash = "18696a239d29cca0989480ad0ccf613ceb5ae8b4"
$newHash = "366bb1c57808d300d691e46ee4af2550ab0c8818"

# the hyperlink target (address) itself is a pinned git-blob URL that is not
# updated by this change -- only the displayed text changes.
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/085fc7af3b4306c8fda636c91b48139bd7c72399/e2e/$oldGuid.md"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $linkAddress, "", "", "e2e\$newGuid.md")
$wsOverview.Range("G2").Value = "2016-08-25 15:00:55"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $linkAddress, "", "", "$newGuid.md")
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-25 15:00:42"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $linkAddress, "", "", "$newGuid.md")
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-25 15:00:55"
